$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest scraped values.
# D-column values are forced to Text via a temporary "@" number format so that
# purely-numeric-looking strings (e.g. "579.48") are not silently turned into
# floating point numbers by Excel's auto-detection; the style is then reset
# back to "Normal" so no stray number-format style is left on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.978.01"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.28%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.114.49"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.67%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.78%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.77%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.110.31"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.62%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.45"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.02%  "

$ws.Range("E11").Value = "  +1.91%  "

$ws.Range("E12").Value = "  +0.89%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000250"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.02%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.30"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.39%  "

$ws.Range("E15").Value = "  +0.25%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.627.41"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.52%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.960.81"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.32%  "

$ws.Range("E18").Value = "  -0.15%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.113.40"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.51%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.45%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "487.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.89%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.719"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.94%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.59"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.66%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.58"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.63%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.63%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.27%  "

$ws.Range("E27").Value = "  +0.38%  "

$ws.Range("E28").Value = "  +0.09%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.06"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.36%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.44"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.22%  "

$ws.Range("E31").Value = "  +2.65%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "29.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.34%  "

$ws.Range("E33").Value = "  +1.79%  "

$ws.Range("E34").Value = "  -3.46%  "

$ws.Range("E35").Value = "  -0.04%  "

$ws.Range("E36").Value = "  +2.02%  "

$ws.Range("E37").Value = "  -0.27%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "47.52"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.97%  "

$ws.Range("E39").Value = "  +3.70%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.12"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.22%  "

$ws.Range("E41").Value = "  +2.53%  "

$ws.Range("E42").Value = "  +1.61%  "

$ws.Range("E43").Value = "  +1.00%  "

$ws.Range("E44").Value = "  -0.97%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.847.19"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.13%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "385.77"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.77%  "

$ws.Range("E47").Value = "  +0.36%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "136.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.88%  "

$ws.Range("E49").Value = "  -0.01%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.20"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.59%  "

$ws.Range("E51").Value = "  +0.69%  "
